$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =========================================================================
# 1) Header row (row 1) text updates
# =========================================================================
$ws.Range("A1").Value2 = "Start Time"
$ws.Range("B1").Value2 = "End Time"
$ws.Range("C1").Value2 = "Duration (hrs)"
$ws.Range("D1").Value2 = "Cost (€)"
$ws.Range("E1").Value2 = "Note"

# =========================================================================
# 2) Pre-touch row 3 with placeholder values so the cells/row formally
#    exist in the sheet's used-range *before* we start copying formats
#    into them. (Formatting a never-before-valued cell confuses the
#    formula dependency graph and leaves SUM() results stale.)
# =========================================================================
$ws.Range("A3:E3").Value2 = 0

# =========================================================================
# 3) Build row 3 by cloning the formatting of the matching column in row 2
#    (keeps each column's own number format: date / hrs / currency / text)
#    E3 (a text/Note cell) clones from G1, which already carries the
#    "no number format, bordered" look that text cells use.
# =========================================================================
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("G1").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# =========================================================================
# 4) Row 2 data: contents rotate C2 -> D2 -> E2 (column formats/styles stay
#    put; only the values/text that live in them change)
# =========================================================================
$ws.Range("C2").Value2 = 1.666666666666667
$ws.Range("D2").Value2 = 41.66666666666667
$ws.Range("E2").Value2 = "Data Cleaning, ho cercato di far quadrare l'excel"

# =========================================================================
# 5) New row 3 values
# =========================================================================
$ws.Range("A3").Value2 = 46061.56944444445
$ws.Range("B3").Value2 = 46061.59930555556
$ws.Range("C3").Value2 = 0.7166666666666667
$ws.Range("D3").Value2 = 17.91666666666667
$ws.Range("E3").Value2 = "Refactored the cleaning code"

# =========================================================================
# 6) Formulas now span rows 2:3
# =========================================================================
$ws.Range("H2").Formula = "=SUM(C2:C3)"
$ws.Range("I2").Formula = "=SUM(D2:D3)"
$excel.Calculate()

# =========================================================================
# 7) Alignment updates
# =========================================================================
# Data columns A:D (dates, duration, cost) -> vertical-center
$ws.Range("A2:D3").VerticalAlignment = -4108

# Note column E -> vertical-center + wrap text
$ws.Range("E2:E3").VerticalAlignment = -4108
$ws.Range("E2:E3").WrapText = $true

# "TOTAL DUE:" label -> right align + vertical center
$ws.Range("G2").HorizontalAlignment = -4152
$ws.Range("G2").VerticalAlignment = -4108

# Totals H2:I2 -> center align + vertical center
$ws.Range("H2:I2").HorizontalAlignment = -4108
$ws.Range("H2:I2").VerticalAlignment = -4108

$excel.Calculate()

# =========================================================================
# 8) Column widths: C,D narrower (15->10); new E column width 30
#    (Excel pads ColumnWidth by 5/6 character on round-trip, so we
#    compensate by subtracting 5/6 from the desired width.)
# =========================================================================
$ws.Columns.Item(3).ColumnWidth = (10 - 0.8333333333333333)
$ws.Columns.Item(4).ColumnWidth = (10 - 0.8333333333333333)
$ws.Columns.Item(5).ColumnWidth = (30 - 0.8333333333333333)
